# Insert a new "RepositoryTest updated and enhanced" bullet item right
# after the "ScoreBoardTest updated to accommodate the changes in
# ScoreBoard" bullet, and before the (empty, bookmarked) trailing
# paragraph of that list / table cell.
#
# We build the new paragraph via raw OOXML (Range.InsertXML) so the
# paragraph properties (ListParagraph style, the same numbering as the
# rest of the list, and the exact spacing) and run properties (Consolas
# font, color, size, highlight) match the list precisely -- something
# that is hard to reproduce reliably purely through the scalar
# Paragraph/Font property setters.

$d = $word.ActiveDocument

# Locate the anchor paragraph: the one whose text is exactly
# "ScoreBoardTest updated to accommodate the changes in ScoreBoard".
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    $text = $candidate.Range.Text
    if ($text -like "*ScoreBoardTest updated to accommodate the changes in ScoreBoard*") {
        $anchor = $candidate
    }
}

$anchorRange = $anchor.Range

# Collapse to a point just before the anchor paragraph's own paragraph
# mark (Range.End - 1), so the freshly-inserted paragraph mark lands
# *before* that mark -- i.e. the new paragraph becomes the paragraph
# right after the anchor, rather than after the following (bookmarked)
# paragraph.
$insertAt = $d.Range($anchorRange.End - 1, $anchorRange.End - 1)

$newParagraphXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr><w:spacing w:before="0" w:after="200" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:noProof/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="2B91AF"/><w:sz w:val="19"/><w:szCs w:val="19"/><w:highlight w:val="white"/></w:rPr><w:t>RepositoryTest</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/><w:color w:val="2B91AF"/><w:sz w:val="19"/><w:szCs w:val="19"/></w:rPr><w:t xml:space="preserve"> updated and enhanced</w:t></w:r></w:p>'

$flatOpc = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertAt.InsertXML($flatOpc)
